# "preice with by items"
# Rework the invoice header block: the big bold/underlined title (A2) becomes a
# normal-sized header line, the address lines are edited/reordered, a new
# "Nuit:" (tax number) line is added, and the old "Itens"/"Preço" two-column
# header row is removed (replaced by a blank spacer row) since the price
# column now lives with the item rows instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. A2 ("Girassol Lavandaria") loses its special bold/underline/12pt look and
#    becomes a normal 8pt line like the rows below it; copy the format from A3
#    (which already uses that style) instead of re-creating fonts by hand, then
#    auto-fit row 2 so its custom 15.75 row height reverts back to the default.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(2).AutoFit() | Out-Null

# 2. Text edits to the address/contact lines.
$ws.Range("A3").Value = "AV. Ahmed Sekou/touré,"
$ws.Range("A6").Value = "cell:+258860852222"

# 3. New line with the company tax id, right after the phone/cell line, using
#    the same formatting as the line above it.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Value = "Nuit:401376070"

# 4. Drop the old "Itens"/"Preço" header cells on row 15: A15 is removed
#    entirely, and B15 becomes a blank spacer cell like B16:B22 below it.
$ws.Range("A15").Clear() | Out-Null
$ws.Range("B16").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("B15").ClearContents() | Out-Null

# 5. Remove one spacer row from the items block (old row 23), shifting the
#    trailing total/footer rows up by one.
$ws.Rows.Item(23).Delete(-4162) | Out-Null

# 6. Widen column A to fit the new longer header lines.
$ws.Columns.Item(1).ColumnWidth = 26.66

# 7. Leave the selection on the new top header cell.
$ws.Range("A2").Select() | Out-Null
